# [TEX] Kupferpaneel; Sprachdatei de geupdated
#
# Adds a second language-key block (columns E:G) mirroring the existing
# A:C block ("item.industria. + <key> + ":) for the newly added copper
# panel related blocks/plates (iron, copper, steel, tin, brass - each in
# exposed/weathered/oxidized and waxed variants).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New language keys to add, in the exact order they must appear
# (row 2 downwards), mirrored into column F while E/G repeat the
# surrounding literal text used by every other row in the sheet.
$newKeys = @(
    "exposed_iron_block",
    "weathered_iron_block",
    "oxidized_iron_block",
    "exposed_iron_plates",
    "weathered_iron_plates",
    "oxidized_iron_plates",
    "waxed_exposed_iron_block",
    "waxed_weathered_iron_block",
    "waxed_oxidized_iron_block",
    "waxed_exposed_iron_plates",
    "waxed_weathered_iron_plates",
    "waxed_oxidized_iron_plates",
    "exposed_copper_plates",
    "weathered_copper_plates",
    "oxidized_copper_plates",
    "waxed_exposed_copper_plates",
    "waxed_weathered_copper_plates",
    "waxed_oxidized_copper_plates",
    "exposed_steel_block",
    "weathered_steel_block",
    "oxidized_steel_block",
    "exposed_steel_plates",
    "weathered_steel_plates",
    "oxidized_steel_plates",
    "waxed_exposed_steel_block",
    "waxed_weathered_steel_block",
    "waxed_oxidized_steel_block",
    "waxed_exposed_steel_plates",
    "waxed_weathered_steel_plates",
    "waxed_oxidized_steel_plates",
    "exposed_tin_block",
    "weathered_tin_block",
    "oxidized_tin_block",
    "exposed_tin_plates",
    "weathered_tin_plates",
    "oxidized_tin_plates",
    "waxed_exposed_tin_block",
    "waxed_weathered_tin_block",
    "waxed_oxidized_tin_block",
    "waxed_exposed_tin_plates",
    "waxed_weathered_tin_plates",
    "waxed_oxidized_tin_plates",
    "exposed_brass_block",
    "weathered_brass_block",
    "oxidized_brass_block",
    "exposed_brass_plates",
    "weathered_brass_plates",
    "oxidized_brass_plates",
    "waxed_exposed_brass_block",
    "waxed_weathered_brass_block",
    "waxed_oxidized_brass_block",
    "waxed_exposed_brass_plates",
    "waxed_weathered_brass_plates",
    "waxed_oxidized_brass_plates"
)

$prefix = '"item.industria.'
$suffix = '":'

$startRow = 2
for ($i = 0; $i -lt $newKeys.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $prefix     # column E
    $ws.Cells.Item($row, 6).Value = $newKeys[$i] # column F
    $ws.Cells.Item($row, 7).Value = $suffix      # column G
}

# Update the active selection to match what was recorded when the
# workbook was saved after this edit.
$ws.Range("I6").Select()
